$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F: header "SBO_DEF" and values "[]" for each data row,
# mirroring the existing header / body styling used by the other columns.
$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("F1").Style = $ws.Range("B1").Style

$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("F4").Value = "[]"
